# Generate Report for Handback
#
# - Marks the localization status as handed back (was "Ready for handoff")
#   everywhere it is shown: Overview!E2:F3 and the "Status" column (C2:C3)
#   on the zh-cn / de-de detail sheets.
# - Stamps the "Latest Handback DateTime" for both languages, and fills in
#   "Latest Target File" / "Latest Handback File" with the generated
#   xliff / markdown info, wiring a hyperlink onto the new Target File
#   cell just like the existing Source File Name hyperlinks.
# - Widens the columns that now hold the longer strings.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackTime = "2016-08-30 08:42:46"
$deHandbackTime = "2016-08-30 08:42:53"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a7e238769f142fbf5cdbb909dbedf5a973dc55c/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: both language status columns flip to "handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = "a.md"
$zh.Range("I3").Value = "a.md"
$zh.Hyperlinks.Add($zh.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null
$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = 15570276
$zh.Range("I3").Font.Underline = 2
$zh.Range("I3").Font.Color = 15570276

$zh.Range("J2").Value = $zhXlf
$zh.Range("J3").Value = $zhXlf

$zh.Range("K2").Value = $zhHandbackTime
$zh.Range("K3").Value = $zhHandbackTime

$zh.Columns.Item(3).ColumnWidth = 29.14
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = "a.md"
$de.Range("I3").Value = "a.md"
$de.Hyperlinks.Add($de.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null
$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = 15570276
$de.Range("I3").Font.Underline = 2
$de.Range("I3").Font.Color = 15570276

$de.Range("J2").Value = $deXlf
$de.Range("J3").Value = $deXlf

$de.Range("K2").Value = $deHandbackTime
$de.Range("K3").Value = $deHandbackTime

$de.Columns.Item(3).ColumnWidth = 29.14
$de.Columns.Item(10).ColumnWidth = 39.17
